# Fixed some color coding in presentation
#
# Four colour-coded "highlight" rectangles get resized/repositioned and, on
# two slides, a new rectangle (colour 3E32EB) is added right after the
# existing one by duplicating it and then restyling the duplicate.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 3 : shape 11 (00B0F0 highlight) resized/moved; new 3E32EB rect added
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$s3shp = $s3.Shapes.Item(11)
$s3shp.Left = 180.6774139404297
$s3shp.Top = 93.62725067138672
$s3shp.Width = 105.15685272216797
$s3shp.Height = 30.28362274169922

$s3new = $s3shp.Duplicate()
$s3new.Left = 314.2460021972656
$s3new.Top = 93.62725067138672
$s3new.Width = 76.98047637939453
$s3new.Height = 30.28362274169922
$s3new.Fill.ForeColor.RGB = 0xEB323E
$s3new.Fill.Transparency = 0.34

# ---------------------------------------------------------------------------
# Slide 4 : shape 11 (accent6 highlight) resized/moved; new 3E32EB rect added
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$s4shp = $s4.Shapes.Item(11)
$s4shp.Left = 186.5597686767578
$s4shp.Top = 93.62732696533203
$s4shp.Width = 102.7058334350586
$s4shp.Height = 30.283544540405273

$s4new = $s4shp.Duplicate()
$s4new.Left = 317.18719482421875
$s4new.Top = 93.62732696533203
$s4new.Width = 76.98047637939453
$s4new.Height = 30.28362274169922
$s4new.Fill.ForeColor.RGB = 0xEB323E
$s4new.Fill.Transparency = 0.34

# ---------------------------------------------------------------------------
# Slide 7 : shape 14 (9B1BDB highlight) resized/moved only
# ---------------------------------------------------------------------------
$s7 = $p.Slides.Item(7)
$s7shp = $s7.Shapes.Item(14)
$s7shp.Left = 624.3048095703125
$s7shp.Top = 90.68614196777344
$s7shp.Width = 111.03937530517578
$s7shp.Height = 33.224884033203125

# ---------------------------------------------------------------------------
# Slide 9 : shape 14 (accent6 highlight) resized/moved only
# ---------------------------------------------------------------------------
$s9 = $p.Slides.Item(9)
$s9shp = $s9.Shapes.Item(14)
$s9shp.Left = 409.8919982910156
$s9shp.Top = 98.95614624023438
$s9shp.Width = 100.94244384765625
$s9shp.Height = 24.954803466796875
